# The authored change swaps the two embedded themes of this deck:
#   ppt/theme/theme1.xml (was the "Integral" color scheme, used by the
#   slide master / presentation) and ppt/theme/theme2.xml (was the
#   default "Office Theme" color scheme, used by the notes master) trade
#   places. Net visible effect for the slide master/presentation theme
#   (the only theme this COM host's object model exposes) is that its
#   12 theme colours change from the "Integral" palette to the stock
#   "Office Theme" palette.

function ToRGB($r, $g, $b) {
    # VBA-style RGB(): COM RGB properties are stored as 0x00BBGGRR.
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Office Theme colour scheme (target values), in the fixed
# ThemeColorScheme.Colors() order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$tcs.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink
